# ragans_cole.xlsx: regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals.
#
# The "K" column (column G, header "K") is recalculated for every game
# row on the sheet. Write the newly computed K values back into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 8
    3  = 5
    4  = 8
    5  = 6
    6  = 7
    7  = 9
    8  = 11
    9  = 9
    10 = 5
    11 = 12
    12 = 8
    13 = 3
    14 = 1
    15 = 0
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 4
    21 = 1
    22 = 2
    23 = 2
    24 = 0
    25 = 4
    26 = 0
    27 = 2
    28 = 3
    29 = 2
    30 = 0
    31 = 4
    32 = 6
    33 = 5
    34 = 4
    35 = 2
    36 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
